$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A147").Value = "备忘录模式"
$ws.Range("B147").Value = "备忘录模式提供了一种方式来捕获对象的内部状态，这样在未来可以将该对象恢复到之前的状态。备忘录模式常用于需要实现撤销功能的场景。"
$ws.Range("C147").Value = "Originator：创建一个备忘录并存储其当前的状态。"
$ws.Range("B148").Value = "优点：`n1. 可以在不破坏封装性的前提下捕获和外部化对象的内部状态。`n2. 可以简化原始对象的实现。"
$ws.Range("C148").Value = "Memento：存储Originator的内部状态。"
$ws.Range("C149").Value = "Caretaker：负责存储备忘录，并在需要时恢复对象的状态。"
$ws.Range("B150").Value = "缺点：`n1. 如果每次保存对象的完整状态，可能会产生大量的存储开销。"
$ws.Range("A152").Value = "原型模式"
$ws.Range("B152").Value = "原型模式涉及到创建对象的克隆。此模式涉及实现一个原型接口，该接口有一个克隆方法。这种模式用于减少创建对象的开销，特别是当创建新对象的成本较大时。"
$ws.Range("B153").Value = "优点：`n1. 与直接创建实例相比，克隆通常更为快速。`n2. 通过克隆可以得到与现有实例具有相同的属性的新实例。"
$ws.Range("C153").Value = "Prototype: 一个声明克隆方法的接口。"
$ws.Range("C154").Value = "ConcretePrototype: 实现了Prototype接口的类。"
$ws.Range("B155").Value = "缺点：`n1. 克隆的对象和源对象不是同一个实例。`n2. 深克隆比浅克隆复杂，并可能涉及复杂的逻辑。"
$ws.Range("A157").Value = "访问者模式"
$ws.Range("B157").Value = "访问者模式用于将数据结构与数据操作分离。它允许添加进一步操作到对象的结构中，而无需修改那些对象。在这种模式中，主要是数据结构提供一个接受访问者的机制，访问者则完成具体的操作。"
$ws.Range("C157").Value = "Visitor (访问者接口): 定义一个访问具体元素的接口，它为每一个具体元素都定义了一个访问操作。"
$ws.Range("B158").Value = "优点：`n1. 增加新的操作很容易，只需增加新的访问者即可。`n2. 分离数据结构和数据操作，符合单一职责原则。`n3. 让每一个访问者角色都不是很庞大，各自的编码也相对简单。"
$ws.Range("C158").Value = "ConcreteVisitor (具体访问者): 实现每个由访问者声明的操作。"
$ws.Range("C159").Value = "Element (元素接口): 定义一个接受访问者的操作。"
$ws.Range("B160").Value = "缺点：`n1. 增加新的数据结构比较困难。`n2. 具体元素变更时，可能需要更改抽象访问者角色。"
$ws.Range("C160").Value = "ConcreteElement (具体元素): 实现接受访问者的操作，并将引用传递给访问者的相应操作。"

$null = $ws.Range("C161").Select()

